$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 34, shifting rows 34:39 down to 35:40
$ws.Rows("34:34").Insert()

# Populate the new row 34 with label and value
$ws.Range("A34").Value = "Ave Pressure kpa (air)"
$ws.Range("B34").Value = 64.56

# Update selection / view state to mirror final saved state
$ws.Range("C34").Select()
